# Add the new "Contacts" worksheet after the existing "Companies" sheet and
# populate it with the new contacts table (Title / First Name / Last Name /
# Company), matching the target diff.

$wb = $excel.ActiveWorkbook

$companies = $wb.Worksheets.Item("Companies")

# Insert a new sheet right after the last existing sheet ("Companies"),
# which also makes it the active sheet/tab (mirrors the diff's activeTab="1"
# and the tabSelected flag moving from sheet1 to the new sheet2).
$contacts = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$contacts.Name = "Contacts"

# Header row
$contacts.Range("A1").Value = "Title"
$contacts.Range("B1").Value = "First Name"
$contacts.Range("C1").Value = "Last Name"
$contacts.Range("D1").Value = "Company"

# Data rows
$contacts.Range("A2").Value = "Dr."
$contacts.Range("B2").Value = "Matt"
$contacts.Range("C2").Value = "Williams"
$contacts.Range("D2").Value = "Mediaocean Asia Pvt. Ltd."

$contacts.Range("A3").Value = "Mr."
$contacts.Range("B3").Value = "Mark"
$contacts.Range("C3").Value = "Keller"
$contacts.Range("D3").Value = "Cognizant Technology Solutions Pvt. Ltd."

$contacts.Range("A4").Value = "Mrs."
$contacts.Range("B4").Value = "Anshika"
$contacts.Range("C4").Value = "Khandelwal"
$contacts.Range("D4").Value = "Infosys Pvt. Ltd."

# Match the header styling used on the "Companies" sheet (yellow fill).
$contacts.Range("A1:D1").Interior.Color = $companies.Range("A1").Interior.Color

# Column widths (best achievable given COM ColumnWidth rounding) matching
# the widths on the "Companies" sheet layout.
$contacts.Columns.Item(1).ColumnWidth = 15
$contacts.Columns.Item(2).ColumnWidth = 19
$contacts.Columns.Item(3).ColumnWidth = 22
$contacts.Columns.Item(4).ColumnWidth = 39.6666666666667

# Selection/active-cell state, matching the sheet's saved selection (A2).
$contacts.Range("A2").Select() | Out-Null
